$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Argentina
$ws.Range("F2").Value = 42.4
$ws.Range("G2").Value = 42.5
$ws.Range("H2").Value = 42.7

# Row 3 - Bolivia
$ws.Range("E3").Value = 42.3
$ws.Range("F3").Value = 42.8
$ws.Range("G3").Value = 43.3
$ws.Range("H3").Value = 43.9

# Row 4 - Brazil
$ws.Range("E4").Value = 50.3
$ws.Range("F4").Value = 50.2
$ws.Range("G4").Value = 50.2
$ws.Range("H4").Value = 50.2

# Row 5 - Chile
$ws.Range("E5").Value = 43.1
$ws.Range("F5").Value = 43.1
$ws.Range("G5").Value = 43.1
$ws.Range("H5").Value = 43.2

# Row 6 - Colombia
$ws.Range("E6").Value = 54.4
$ws.Range("F6").Value = 54.2
$ws.Range("G6").Value = 54.1
$ws.Range("H6").Value = 54.3

# Row 7 - Costa Rica
$ws.Range("F7").Value = 45.7
$ws.Range("G7").Value = 45.7
$ws.Range("H7").Value = 45.7

# Row 8 - Dominican Republic
$ws.Range("E8").Value = 39
$ws.Range("F8").Value = 38.8
$ws.Range("G8").Value = 38.8
$ws.Range("H8").Value = 38.8

# Row 9 - Ecuador
$ws.Range("G9").Value = 45.2
$ws.Range("H9").Value = 45.2

# Row 10 - Guatemala
$ws.Range("E10").Value = 45.4
$ws.Range("F10").Value = 45.7
$ws.Range("H10").Value = 46.4

# Row 11 - Honduras
$ws.Range("E11").Value = 45.7
$ws.Range("F11").Value = 45.4
$ws.Range("G11").Value = 45.5
$ws.Range("H11").Value = 45.1

# Row 12 - Latin America and the Caribbean: remove all data values
$ws.Range("C12:H12").ClearContents()

# Row 14 - Nicaragua
$ws.Range("E14").Value = 49.4
$ws.Range("F14").Value = 49.9
$ws.Range("G14").Value = 49.6
$ws.Range("H14").Value = 49.5

# Row 15 - Panama
$ws.Range("E15").Value = 49.7
$ws.Range("F15").Value = 49.7
$ws.Range("G15").Value = 49.8
$ws.Range("H15").Value = 49.7

# Row 16 - Peru
$ws.Range("E16").Value = 40.1
$ws.Range("F16").Value = 40.3
$ws.Range("G16").Value = 40
$ws.Range("H16").Value = 39.8

# Row 17 - Paraguay
$ws.Range("C17").Value = 44.6
$ws.Range("D17").Value = 44.2
$ws.Range("E17").Value = 44.2
$ws.Range("F17").Value = 44.5
$ws.Range("G17").Value = 44.7
$ws.Range("H17").Value = 44.7

# Row 18 - El Salvador
$ws.Range("E18").Value = 40
$ws.Range("F18").Value = 40.2
$ws.Range("G18").Value = 40.4
$ws.Range("H18").Value = 40.5

# Row 19 - Uruguay
$ws.Range("E19").Value = 40
$ws.Range("F19").Value = 40
$ws.Range("G19").Value = 40.1
$ws.Range("H19").Value = 40.1
